$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 1199.125
$ws.Cells.Item(40, 9).Value = 1100.25
$ws.Cells.Item(40, 10).Value = 1298
$ws.Cells.Item(40, 11).Value = 1100.25
$ws.Cells.Item(40, 12).Value = 1298
$ws.Cells.Item(40, 13).Value = -925.25
$ws.Cells.Item(40, 14).Value = -1648
# Row 42
$ws.Cells.Item(42, 8).Value = 457.27777
$ws.Cells.Item(42, 9).Value = 462.75
$ws.Cells.Item(42, 10).Value = 452.9
$ws.Cells.Item(42, 11).Value = 1388.25
$ws.Cells.Item(42, 12).Value = 1358.7
$ws.Cells.Item(42, 13).Value = -1158.25
$ws.Cells.Item(42, 14).Value = -1818.7
# Row 62
$ws.Cells.Item(62, 8).Value = 1455.0834
$ws.Cells.Item(62, 9).Value = 1591
$ws.Cells.Item(62, 10).Value = 1358
$ws.Cells.Item(62, 11).Value = 1591
$ws.Cells.Item(62, 12).Value = 1358
$ws.Cells.Item(62, 13).Value = -967
$ws.Cells.Item(62, 14).Value = -2606
# Row 65
$ws.Cells.Item(65, 8).Value = 1455.0834
$ws.Cells.Item(65, 9).Value = 1591
$ws.Cells.Item(65, 10).Value = 1358
$ws.Cells.Item(65, 11).Value = 7955
$ws.Cells.Item(65, 12).Value = 6790
$ws.Cells.Item(65, 13).Value = -4835
$ws.Cells.Item(65, 14).Value = -13030
# Row 99
$ws.Cells.Item(99, 8).Value = 1250
$ws.Cells.Item(99, 10).Value = 2000
$ws.Cells.Item(99, 12).Value = 6000
$ws.Cells.Item(99, 14).Value = -8996
# Row 132
$ws.Cells.Item(132, 8).Value = 5272.7295
$ws.Cells.Item(132, 9).Value = 4686.4194
$ws.Cells.Item(132, 10).Value = 8302
$ws.Cells.Item(132, 11).Value = 14059.2582
$ws.Cells.Item(132, 12).Value = 24906
$ws.Cells.Item(132, 13).Value = -11529.2582
$ws.Cells.Item(132, 14).Value = -29966

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 2016.174
$ws.Cells.Item(61, 9).Value = 1490.1428
$ws.Cells.Item(61, 11).Value = 1490.1428
$ws.Cells.Item(61, 13).Value = -1278.1428
# Row 122
$ws.Cells.Item(122, 8).Value = 2067.6086
$ws.Cells.Item(122, 9).Value = 1478.8
$ws.Cells.Item(122, 11).Value = 4436.4
$ws.Cells.Item(122, 13).Value = -1986.4
# Row 136
$ws.Cells.Item(136, 8).Value = 2016.174
$ws.Cells.Item(136, 9).Value = 1490.1428
$ws.Cells.Item(136, 11).Value = 4470.428400000001
$ws.Cells.Item(136, 13).Value = -1920.428400000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 1202.6296
$ws.Cells.Item(107, 9).Value = 1203.45
$ws.Cells.Item(107, 10).Value = 1200.2858
$ws.Cells.Item(107, 11).Value = 1203.45
$ws.Cells.Item(107, 12).Value = 1200.2858
$ws.Cells.Item(107, 13).Value = 716.55
$ws.Cells.Item(107, 14).Value = -5040.2858
# Row 122
$ws.Cells.Item(122, 8).Value = 41768.75
$ws.Cells.Item(122, 10).Value = 41768.75
$ws.Cells.Item(122, 12).Value = 41768.75
$ws.Cells.Item(122, 14).Value = -51568.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 16670655
$ws.Cells.Item(31, 9).Value = 2352.158
$ws.Cells.Item(31, 10).Value = 45461360
$ws.Cells.Item(31, 11).Value = 2352.158
$ws.Cells.Item(31, 12).Value = 45461360
$ws.Cells.Item(31, 13).Value = -2057.158
$ws.Cells.Item(31, 14).Value = -45461950
# Row 34
$ws.Cells.Item(34, 8).Value = 16670655
$ws.Cells.Item(34, 9).Value = 2352.158
$ws.Cells.Item(34, 10).Value = 45461360
$ws.Cells.Item(34, 11).Value = 2352.158
$ws.Cells.Item(34, 12).Value = 45461360
$ws.Cells.Item(34, 13).Value = -2150.158
$ws.Cells.Item(34, 14).Value = -45461764
# Row 122
$ws.Cells.Item(122, 8).Value = 1854.4615
$ws.Cells.Item(122, 9).Value = 1242.6666
$ws.Cells.Item(122, 10).Value = 2178.353
$ws.Cells.Item(122, 11).Value = 3727.9998
$ws.Cells.Item(122, 12).Value = 6535.059
$ws.Cells.Item(122, 13).Value = -1277.9998
$ws.Cells.Item(122, 14).Value = -11435.059

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Cells.Item(21, 8).Value = 2646.8262
$ws.Cells.Item(21, 9).Value = 546.5
$ws.Cells.Item(21, 11).Value = 1639.5
$ws.Cells.Item(21, 13).Value = -1466.5
# Row 39
$ws.Cells.Item(39, 8).Value = 11984
$ws.Cells.Item(39, 10).Value = 11984
$ws.Cells.Item(39, 12).Value = 35952
$ws.Cells.Item(39, 14).Value = -36540
# Row 70
$ws.Cells.Item(70, 8).Value = 1627.4
$ws.Cells.Item(70, 9).Value = 1128.2727
$ws.Cells.Item(70, 11).Value = 3384.8181
$ws.Cells.Item(70, 13).Value = -3069.8181
# Row 73
$ws.Cells.Item(73, 8).Value = 1627.4
$ws.Cells.Item(73, 9).Value = 1128.2727
$ws.Cells.Item(73, 11).Value = 3384.8181
$ws.Cells.Item(73, 13).Value = -2292.8181
# Row 109
$ws.Cells.Item(109, 8).Value = 5110.5835
$ws.Cells.Item(109, 9).Value = 756.7143
$ws.Cells.Item(109, 11).Value = 2270.1429
$ws.Cells.Item(109, 13).Value = -1230.1429
# Row 112
$ws.Cells.Item(112, 8).Value = 5123.8096
$ws.Cells.Item(112, 10).Value = 5433.3335
$ws.Cells.Item(112, 12).Value = 16300.0005
$ws.Cells.Item(112, 14).Value = -18516.0005
# Row 113
$ws.Cells.Item(113, 8).Value = 641.7857
$ws.Cells.Item(113, 9).Value = 618
$ws.Cells.Item(113, 10).Value = 653.6786
$ws.Cells.Item(113, 11).Value = 1854
$ws.Cells.Item(113, 12).Value = 1961.0358
$ws.Cells.Item(113, 13).Value = 316
$ws.Cells.Item(113, 14).Value = -6301.0358
# Row 117
$ws.Cells.Item(117, 8).Value = 3506.842
$ws.Cells.Item(117, 10).Value = 3668.3333
$ws.Cells.Item(117, 12).Value = 11004.9999
$ws.Cells.Item(117, 14).Value = -17888.9999
# Row 121
$ws.Cells.Item(121, 8).Value = 1721.8167
$ws.Cells.Item(121, 9).Value = 300
$ws.Cells.Item(121, 10).Value = 1745.9153
$ws.Cells.Item(121, 11).Value = 900
$ws.Cells.Item(121, 12).Value = 5237.7459
$ws.Cells.Item(121, 13).Value = 410
$ws.Cells.Item(121, 14).Value = -7857.7459
# Row 127
$ws.Cells.Item(127, 8).Value = 1143.3334
$ws.Cells.Item(127, 10).Value = 1143.3334
$ws.Cells.Item(127, 12).Value = 3430.0002
$ws.Cells.Item(127, 14).Value = -13350.0002
# Row 131
$ws.Cells.Item(131, 8).Value = 7353792.5
$ws.Cells.Item(131, 10).Value = 860.8594000000001
$ws.Cells.Item(131, 12).Value = 2582.5782
$ws.Cells.Item(131, 14).Value = -12662.5782
# Row 137
$ws.Cells.Item(137, 8).Value = 3741.5
$ws.Cells.Item(137, 9).Value = 500
$ws.Cells.Item(137, 10).Value = 6983
$ws.Cells.Item(137, 11).Value = 1500
$ws.Cells.Item(137, 12).Value = 20949
$ws.Cells.Item(137, 13).Value = 3600
$ws.Cells.Item(137, 14).Value = -31149

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Cells.Item(5, 8).Value = 10991.667
$ws.Cells.Item(5, 10).Value = 10991.667
$ws.Cells.Item(5, 12).Value = 10991.667
$ws.Cells.Item(5, 14).Value = -11215.667
# Row 97
$ws.Cells.Item(97, 8).Value = 758.1667
$ws.Cells.Item(97, 9).Value = 709.8
$ws.Cells.Item(97, 11).Value = 709.8
$ws.Cells.Item(97, 13).Value = -213.8
# Row 102
$ws.Cells.Item(102, 8).Value = 1743.2
$ws.Cells.Item(102, 9).Value = 1302.0741
$ws.Cells.Item(102, 10).Value = 2659.3845
$ws.Cells.Item(102, 11).Value = 1302.0741
$ws.Cells.Item(102, 12).Value = 2659.3845
$ws.Cells.Item(102, 13).Value = 319.9259
$ws.Cells.Item(102, 14).Value = -5903.3845
# Row 122
$ws.Cells.Item(122, 8).Value = 2719.5386
$ws.Cells.Item(122, 9).Value = 1481.375
$ws.Cells.Item(122, 10).Value = 4700.6
$ws.Cells.Item(122, 11).Value = 4444.125
$ws.Cells.Item(122, 12).Value = 14101.8
$ws.Cells.Item(122, 13).Value = -1994.125
$ws.Cells.Item(122, 14).Value = -19001.8
# Row 126
$ws.Cells.Item(126, 8).Value = 2829.7
$ws.Cells.Item(126, 9).Value = 2870.2395
$ws.Cells.Item(126, 10).Value = 1856.75
$ws.Cells.Item(126, 11).Value = 8610.718500000001
$ws.Cells.Item(126, 12).Value = 5570.25
$ws.Cells.Item(126, 13).Value = -6140.718500000001
$ws.Cells.Item(126, 14).Value = -10510.25
# Row 132
$ws.Cells.Item(132, 8).Value = 2773.4546
$ws.Cells.Item(132, 9).Value = 1376.9445
$ws.Cells.Item(132, 10).Value = 4449.2666
$ws.Cells.Item(132, 11).Value = 4130.833500000001
$ws.Cells.Item(132, 12).Value = 13347.7998
$ws.Cells.Item(132, 13).Value = -1600.833500000001
$ws.Cells.Item(132, 14).Value = -18407.7998

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 2690.1
$ws.Cells.Item(46, 9).Value = 2333.5
$ws.Cells.Item(46, 10).Value = 3225
$ws.Cells.Item(46, 11).Value = 2333.5
$ws.Cells.Item(46, 12).Value = 3225
$ws.Cells.Item(46, 13).Value = -2145.5
$ws.Cells.Item(46, 14).Value = -3601
# Row 122
$ws.Cells.Item(122, 8).Value = 4979
$ws.Cells.Item(122, 9).Value = 3076.3572
$ws.Cells.Item(122, 11).Value = 9229.071599999999
$ws.Cells.Item(122, 13).Value = -6779.071599999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 2980.2
$ws.Cells.Item(122, 9).Value = 1921
$ws.Cells.Item(122, 11).Value = 5763
$ws.Cells.Item(122, 13).Value = -3313
